# SetRules.xlsx edit script
# 1. VEDA_Sets-Comm ("Csets"): sort the data rows (3:65) ascending by column D (SetName)
# 2. VEDA_Sets-Proc ("Psets"): append 12 new rows describing residential (RSD) energy-service sets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sort VEDA_Sets-Comm data rows ascending by column D
# ---------------------------------------------------------------------------
$wsComm = $wb.Worksheets.Item("VEDA_Sets-Comm")

$sortRange = $wsComm.Range("A3:J65")
$sortKey = $wsComm.Range("D3:D65")

$wsComm.Sort.SortFields.Clear()
$wsComm.Sort.SortFields.Add($sortKey, 0, 1, 0, 0) | Out-Null
$wsComm.Sort.SetRange($sortRange)
$wsComm.Sort.Header = 0
$wsComm.Sort.Orientation = 1
$wsComm.Sort.Apply()

# ---------------------------------------------------------------------------
# 2) Add the new RSD_NRGSRV-* rows to VEDA_Sets-Proc
# ---------------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("VEDA_Sets-Proc")

$newRows = @(
    @{ Row = 37; B = "R-RSDCD*"; F = "RSD_NRGSRV-CD"; G = "Residential Cloth Drying" },
    @{ Row = 38; B = "R-RSDCK*"; F = "RSD_NRGSRV-CK"; G = "Residential Cooking" },
    @{ Row = 39; B = "R-RSDCW*"; F = "RSD_NRGSRV-CW"; G = "Residential Cloth Washing" },
    @{ Row = 40; B = "R-RSDDW*"; F = "RSD_NRGSRV-DW"; G = "Residential Dish Washing" },
    @{ Row = 41; B = "R-LT*";    F = "RSD_NRGSRV-LT"; G = "Residential Lighting" },
    @{ Row = 42; B = "R-RSDOA*"; F = "RSD_NRGSRV-OA"; G = "Residential Other Applications" },
    @{ Row = 43; B = "R-RSDOE*"; F = "RSD_NRGSRV-OE"; G = "Residential Electric Appliances" },
    @{ Row = 44; B = "R-PF*";    F = "RSD_NRGSRV-PF"; G = "Residential Pump and Fans" },
    @{ Row = 45; B = "R-RSDRF*"; F = "RSD_NRGSRV-RF"; G = "Residential Refrigeration" },
    @{ Row = 46; B = "R-SC*";    F = "RSD_NRGSRV-SC"; G = "Residential Space Cooling" },
    @{ Row = 47; B = "R-SH*";    F = "RSD_NRGSRV-SH"; G = "Residential Space Heating" },
    @{ Row = 48; B = "R-WH*";    F = "RSD_NRGSRV-WH"; G = "Residential Water Heating" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $wsProc.Cells.Item($r, 2).Value2 = $item.B   # B
    $wsProc.Cells.Item($r, 6).Value2 = $item.F   # F
    $wsProc.Cells.Item($r, 7).Value2 = $item.G   # G
    $wsProc.Cells.Item($r, 8).Value2 = "AND"     # H
    $wsProc.Cells.Item($r, 9).Value2 = "OR"      # I
    $wsProc.Cells.Item($r, 10).Value2 = "AND"    # J
    $wsProc.Cells.Item($r, 11).Value2 = "OR"     # K
}

Write-Host "Edit complete"
